$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marcus")

# Update the reported hours for H4 (causes K4, H12, C15 to recalc via formulas)
$ws.Range("H4").Value = 16

# Update the active cell selection to H5
$ws.Range("H5").Select()
